# The <id>...</id> tag pair for p152v_1 was previously split across three
# separate runs:
#   run1 (Courier New, color 7f6000, sz 18): "<id>"
#   run2 (default font, color 000000):       "p152v_1"
#   run3 (Courier New, color 7f6000, sz 18): "</id>"
# Collapse them into a single run (keeping run1's formatting, as Word does
# automatically when the run boundaries disappear) whose text is the full
# concatenation "<id>p152v_1</id>".

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("<id>p152v_1</id>", $true, $false, $false, $false, `
                            $false, $true, 1, $false, "<id>p152v_1</id>", 2)

if (-not $found) {
    throw "Could not find the target '<id>p152v_1</id>' text to update."
}
